$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step A: create new merged ranges first (Merge() resets the anchor cell's
#     style, so do this before writing values/formats into those cells) ---
$ws.Range("C14:C19").Merge()
$ws.Range("D26:D31").Merge()

# --- Step B: remove the old C34:C39 merge (contents/format are cleared later) ---
$ws.Range("C34:C39").UnMerge()

# --- Step C: set new text content (course moved / room changed) ---
$ws.Range("C6").Value = "Scientific Inquiry: Beyond the Visible`n09:00-10:30`nroom:Red classroom: 201"
$ws.Range("E6").Value = "Academic Writing: Research, Fiction and Nonfiction`n09:00-10:30`nroom:Green classroom: 204"

$ws.Range("B14").Value = "Academic Writing: Research, Fiction and Nonfiction`n11:00-12:30`nroom:Green classroom: 204"
$ws.Range("C14").Value = "Academic Writing: Research, Fiction and Nonfiction`n11:00-12:30`nroom:Grey classroom: 203"
$ws.Range("D14").Value = "Precalculus`n11:00-12:30`nroom:Green classroom: 204"

$ws.Range("D26").Value = "Precalculus`n14:00-15:30`nroom:Green classroom: 204"
$ws.Range("E26").Value = "Academic Writing: Research, Fiction and Nonfiction`n14:00-15:30`nroom:Grey classroom: 203"

$ws.Range("E34").Value = "Academic Writing: Research, Fiction and Nonfiction`n16:00-17:30`nroom:Grey classroom: 203"

# --- Step D: copy cell formatting (fill colour / border / font) from cells
#     that already carry the desired look, so the style table is reused
#     instead of growing with near-duplicate entries ---

# red style (used previously only by F14)
$ws.Range("F14").Copy()
$ws.Range("C6").PasteSpecial(-4122)

# green style (used previously by B6/D6)
$ws.Range("B6").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("D26").PasteSpecial(-4122)

# grey style (used previously by C26/B34)
$ws.Range("C26").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("E34").PasteSpecial(-4122)

# --- Step E: give the new blank filler cells (under the new merges) the
#     same "empty slot" borders as their neighbouring columns ---
$ws.Range("B7").Copy()
$ws.Range("C15:C18").PasteSpecial(-4122)
$ws.Range("B11").Copy()
$ws.Range("C19").PasteSpecial(-4122)

$ws.Range("C27").Copy()
$ws.Range("D27:D30").PasteSpecial(-4122)
$ws.Range("C31").Copy()
$ws.Range("D31").PasteSpecial(-4122)

# --- Step F: remove the old C34:C39 column content entirely ---
$ws.Range("C34:C39").Clear()

Write-Host "edit complete"
